$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B11 currently holds the text "R40" (row "Rule" label). It needs to become
# the literal text "1" -- a string, not the number 1. A plain
#   $ws.Range("B11").Value = "1"
# would be auto-coerced to the numeric value 1 by Excel's normal cell-entry
# rules (exactly like typing 1 into a General-formatted cell), which is not
# what we want here (the authoritative edit keeps this a text/shared-string
# cell). To force text entry while leaving B11's own formatting untouched,
# stage the literal text value in a scratch cell (quote-prefixed so Excel
# treats it as text) and copy only the *value* over with PasteSpecial -
# this carries the text without touching B11's existing style/number format.
$scratch = $ws.Range("Z1")
$scratch.Value = "'1"
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)   # xlPasteValues
$scratch.Clear()
